# Tendor_BOM.xlsx: drop the duplicate "ReqBalanceQty" column (D) and keep the
# "TendorQuantity" column, fixing its header typo to "TenderQuantity".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the whole column D ("ReqBalanceQty") - everything to the right
# (old TendorQuantity/ItemUOM/ItemRate/ItemDescription) shifts left.
$ws.Range("D1").EntireColumn.Delete()

# Fix the header typo: "TendorQuantity" -> "TenderQuantity" (now in column D).
$ws.Range("D1").Value = "TenderQuantity"

# Match the author's final selection.
$ws.Range("D2").Select()
